$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Correspond Handback DateTime (Overview) / Correspond Handoff Datetime (de-de)
$wsOverview.Range("G2").Value = "2016-11-02 05:17:00"
$wsOverview.Range("G3").Value = "2016-11-02 05:17:00"
$wsDeDe.Range("H2").Value = "2016-11-02 05:17:00"
$wsDeDe.Range("H3").Value = "2016-11-02 05:17:00"

# Priority: ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# zh-cn Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-11-02 05:16:45"
$wsZhCn.Range("H3").Value = "2016-11-02 05:16:45"

# zh-cn Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-11-02 05:17:49"
$wsZhCn.Range("K3").Value = "2016-11-02 05:17:49"

# de-de Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-11-02 05:18:08"
$wsDeDe.Range("K3").Value = "2016-11-02 05:18:08"
